$wb = $excel.ActiveWorkbook

# Sheet ALC, row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2998.3333
$ws.Range("I62").Value = 2998.3333
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2998.3333
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2374.3333
$ws.Range("N62").ClearContents()

# Sheet ALC, row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 2998.3333
$ws.Range("I65").Value = 2998.3333
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 14991.6665
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -11871.6665
$ws.Range("N65").ClearContents()

# Sheet ALC, row 106
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 3832.3333
$ws.Range("I106").Value = 3832.3333
$ws.Range("K106").Value = 3832.3333
$ws.Range("M106").Value = -3201.3333

# Sheet ALC, row 114
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

# Sheet ALC, row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1438.1666
$ws.Range("I135").Value = 1438.1666
$ws.Range("K135").Value = 12943.4994
$ws.Range("M135").Value = -10408.4994

# Sheet ARM, row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5112.25
$ws.Range("I2").Value = 5112.25
$ws.Range("K2").Value = 5112.25
$ws.Range("M2").Value = -4999.25

# Sheet ARM, row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2393.2666
$ws.Range("I32").Value = 2393.2666
$ws.Range("K32").Value = 2393.2666
$ws.Range("M32").Value = -2106.2666

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2968.9333
$ws.Range("I61").Value = 3003
$ws.Range("J61").Value = 2832.6667
$ws.Range("K61").Value = 3003
$ws.Range("L61").Value = 2832.6667
$ws.Range("M61").Value = -2791
$ws.Range("N61").Value = -3256.6667

# Sheet ARM, row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 8999.666999999999
$ws.Range("I74").Value = 5999.5
$ws.Range("K74").Value = 5999.5
$ws.Range("M74").Value = -5125.5

# Sheet ARM, row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 8999.666999999999
$ws.Range("I77").Value = 5999.5
$ws.Range("K77").Value = 29997.5
$ws.Range("M77").Value = -25629.5

# Sheet ARM, row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 5112.25
$ws.Range("I116").Value = 5112.25
$ws.Range("K116").Value = 5112.25
$ws.Range("M116").Value = -2818.25

# Sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2968.9333
$ws.Range("I136").Value = 3003
$ws.Range("J136").Value = 2832.6667
$ws.Range("K136").Value = 9009
$ws.Range("L136").Value = 8498.000100000001
$ws.Range("M136").Value = -6459
$ws.Range("N136").Value = -13598.0001

# Sheet BSM, row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5112.25
$ws.Range("I3").Value = 5112.25
$ws.Range("K3").Value = 5112.25
$ws.Range("M3").Value = -4998.25

# Sheet BSM, row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6321.3
$ws.Range("I134").Value = 844.2
$ws.Range("K134").Value = 2532.6
$ws.Range("M134").Value = 2.399999999999636

# Sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6675.143
$ws.Range("I31").Value = 6823.6665
$ws.Range("J31").Value = 6563.75
$ws.Range("K31").Value = 6823.6665
$ws.Range("L31").Value = 6563.75
$ws.Range("M31").Value = -6528.6665
$ws.Range("N31").Value = -7153.75

# Sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 6675.143
$ws.Range("I34").Value = 6823.6665
$ws.Range("J34").Value = 6563.75
$ws.Range("K34").Value = 6823.6665
$ws.Range("L34").Value = 6563.75
$ws.Range("M34").Value = -6621.6665
$ws.Range("N34").Value = -6967.75

# Sheet CRP, row 50
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 23054.363
$ws.Range("I50").Value = 15000
$ws.Range("J50").Value = 24844.223
$ws.Range("K50").Value = 15000
$ws.Range("L50").Value = 24844.223
$ws.Range("M50").Value = -14375
$ws.Range("N50").Value = -26094.223

# Sheet CRP, row 51
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 21416.666
$ws.Range("J51").Value = 21416.666
$ws.Range("L51").Value = 21416.666
$ws.Range("N51").Value = -22888.666

# Sheet CRP, row 53
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H53").Value = 55000
$ws.Range("J53").Value = 55000
$ws.Range("L53").Value = 55000
$ws.Range("N53").Value = -56214

# Sheet CRP, row 60
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 19178.572
$ws.Range("J60").Value = 25250
$ws.Range("L60").Value = 25250
$ws.Range("N60").Value = -26272

# Sheet CRP, row 61
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 21416.666
$ws.Range("J61").Value = 21416.666
$ws.Range("L61").Value = 21416.666
$ws.Range("N61").Value = -22112.666

# Sheet CUL, row 3
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3300
$ws.Range("J3").Value = 3300
$ws.Range("L3").Value = 9900
$ws.Range("N3").Value = -10124

# Sheet CUL, row 11
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 125000420
$ws.Range("I11").Value = 166667060
$ws.Range("K11").Value = 500001180
$ws.Range("M11").Value = -500001040

# Sheet CUL, row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1275.4
$ws.Range("J68").Value = 1263.5
$ws.Range("L68").Value = 3790.5
$ws.Range("N68").Value = -5412.5

# Sheet CUL, row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1275.4
$ws.Range("J71").Value = 1263.5
$ws.Range("L71").Value = 11371.5
$ws.Range("N71").Value = -19483.5

# Sheet CUL, row 86
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 1450
$ws.Range("I86").Value = 900
$ws.Range("K86").Value = 2700
$ws.Range("M86").Value = -1514

# Sheet CUL, row 89
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 1450
$ws.Range("I89").Value = 900
$ws.Range("K89").Value = 8100
$ws.Range("M89").Value = -2172

# Sheet GSM, row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2040.1111
$ws.Range("I132").Value = 1694.6
$ws.Range("J132").Value = 2472
$ws.Range("K132").Value = 5083.799999999999
$ws.Range("L132").Value = 7416
$ws.Range("M132").Value = -2553.799999999999
$ws.Range("N132").Value = -12476

# Sheet LTW, row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4858.1665
$ws.Range("I132").Value = 4866.4443
$ws.Range("K132").Value = 14599.3329
$ws.Range("M132").Value = -12069.3329

# Sheet LTW, row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4849.875
$ws.Range("I136").Value = 4849.875
$ws.Range("K136").Value = 14549.625
$ws.Range("M136").Value = -11999.625

# Sheet WVR, row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3961.4443
$ws.Range("J62").Value = 4108.8335
$ws.Range("L62").Value = 4108.8335
$ws.Range("N62").Value = -5356.8335

# Sheet WVR, row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 3961.4443
$ws.Range("J65").Value = 4108.8335
$ws.Range("L65").Value = 20544.1675
$ws.Range("N65").Value = -26784.1675

# Sheet WVR, row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4953
$ws.Range("I122").Value = 3691.25
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 11073.75
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -8623.75
$ws.Range("N122").Value = -34900

# Sheet WVR, row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1970.4286
$ws.Range("I126").Value = 1970.4286
$ws.Range("K126").Value = 5911.2858
$ws.Range("M126").Value = -3441.2858

# Sheet WVR, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4205.4375
$ws.Range("I132").Value = 3295.1538
$ws.Range("K132").Value = 9885.4614
$ws.Range("M132").Value = -7355.4614
